# Correcciones generales, actualización de vistas y datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a genuinely empty TEXT value into a cell (as opposed to
# clearing it to a truly blank cell). A bare "Value = ''' via COM collapses
# the cell to blank, so we go through the classic "force text" leading
# apostrophe, then strip the resulting quotePrefix style back to Normal so
# no stray formatting is left behind.
function Set-EmptyText($rng) {
    $rng.Value = "'"
    $rng.Style = "Normal"
}

# --- Row 5 (veh_1765502675870_3_d333phqbv) ---------------------------------
Set-EmptyText $ws.Range("C5")
Set-EmptyText $ws.Range("G5")
Set-EmptyText $ws.Range("H5")
Set-EmptyText $ws.Range("I5")
$ws.Range("J5").Value = $true

# --- Row 6 (veh_1765502675870_4_3olg50y2q) ---------------------------------
$ws.Range("I6").Value = "Hola"

# --- favorito (J) flips from FALSE to TRUE on several rows ------------------
$ws.Range("J10").Value = $true
$ws.Range("J12").Value = $true
$ws.Range("J13").Value = $true
$ws.Range("J14").Value = $true
$ws.Range("J16").Value = $true

# --- Row 18 (veh_1765502675870_16_8wk4eqtxu): fix brand + mark favorito ----
$ws.Range("D18").Value = "Great Wall"
$ws.Range("J18").Value = $true

# --- New row 19: new vehicle ------------------------------------------------
$ws.Range("A19").Value = "veh_1766019874075_3802"
$ws.Range("C19").Value = "camioneta"
$ws.Range("D19").Value = "Toyota"
$ws.Range("E19").Value = "Hilux 2.4"
$ws.Range("F19").Value = "dasdas"
$ws.Range("G19").Value = "dasd"
$ws.Range("H19").Value = "sadas"
Set-EmptyText $ws.Range("I19")
